# Update countries & provincias Spain
# - Re-rank a handful of countries whose case counts were refreshed
#   (their row keeps its position, but since the table is sorted by
#   total cases descending, the country name + stats for that row change).
# - Refresh the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Banglades (row 19) - updated counts, same rank
$ws.Cells.Item(19, 2).Value = 367565
$ws.Cells.Item(19, 3).Value = 1182
$ws.Cells.Item(19, 4).Value = 280069
$ws.Cells.Item(19, 5).Value = 82171
$ws.Cells.Item(19, 7).Value = 20
$ws.Cells.Item(19, 8).Value = 5325

# Indonesia overtakes Alemania (rows 25-26 swap rank)
$ws.Cells.Item(25, 1).Value = "Indonesia"
$ws.Cells.Item(25, 2).Value = 299506
$ws.Cells.Item(25, 3).Value = 4007
$ws.Cells.Item(25, 4).Value = 225052
$ws.Cells.Item(25, 5).Value = 63399
$ws.Cells.Item(25, 7).Value = 83
$ws.Cells.Item(25, 8).Value = 11055

$ws.Cells.Item(26, 1).Value = "Alemania"
$ws.Cells.Item(26, 2).Value = 298363
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 259500
$ws.Cells.Item(26, 5).Value = 29267
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 9596

# Suiza (row 62) - updated counts, same rank
$ws.Cells.Item(62, 5).Value = 6508
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = 2076

# Austria overtakes Kirguistan and Ghana (rows 65-67 rotate rank)
$ws.Cells.Item(65, 1).Value = "Austria"
$ws.Cells.Item(65, 2).Value = 47432
$ws.Cells.Item(65, 3).Value = 1058
$ws.Cells.Item(65, 4).Value = 38045
$ws.Cells.Item(65, 5).Value = 8578
$ws.Cells.Item(65, 7).Value = 6
$ws.Cells.Item(65, 8).Value = 809

$ws.Cells.Item(66, 1).Value = "Kirguistan"
$ws.Cells.Item(66, 2).Value = 47056
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 43137
$ws.Cells.Item(66, 5).Value = 2854
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 1065

$ws.Cells.Item(67, 1).Value = "Ghana"
$ws.Cells.Item(67, 2).Value = 46694
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 4).Value = 45945
$ws.Cells.Item(67, 5).Value = 448
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 301

# El Salvador (row 78) - updated counts, same rank
$ws.Cells.Item(78, 2).Value = 29358
$ws.Cells.Item(78, 3).Value = 183
$ws.Cells.Item(78, 4).Value = 24175
$ws.Cells.Item(78, 5).Value = 4326
$ws.Cells.Item(78, 7).Value = 4
$ws.Cells.Item(78, 8).Value = 857

# Malasia (row 99) - updated counts, same rank
$ws.Cells.Item(99, 2).Value = 12088
$ws.Cells.Item(99, 3).Value = 317
$ws.Cells.Item(99, 4).Value = 10216
$ws.Cells.Item(99, 5).Value = 1735
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 137

# Finlandia (row 105) - updated counts, same rank
$ws.Cells.Item(105, 2).Value = 10391
$ws.Cells.Item(105, 3).Value = 147
$ws.Cells.Item(105, 5).Value = 1946

# Eslovenia overtakes Cabo Verde (rows 117-118 swap rank)
$ws.Cells.Item(117, 1).Value = "Eslovenia"
$ws.Cells.Item(117, 2).Value = 6330
$ws.Cells.Item(117, 3).Value = 227
$ws.Cells.Item(117, 4).Value = 4041
$ws.Cells.Item(117, 5).Value = 2134
$ws.Cells.Item(117, 7).Value = 1
$ws.Cells.Item(117, 8).Value = 155

$ws.Cells.Item(118, 1).Value = "Cabo Verde"
$ws.Cells.Item(118, 2).Value = 6205
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = 5399
$ws.Cells.Item(118, 5).Value = 744
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 62

# Hong Kong (row 126) - updated counts, same rank
$ws.Cells.Item(126, 2).Value = 5109
$ws.Cells.Item(126, 3).Value = 4
$ws.Cells.Item(126, 4).Value = 4849
$ws.Cells.Item(126, 5).Value = 155

# Lituania overtakes Guinea Ecuatorial (rows 128-129 swap rank)
$ws.Cells.Item(128, 1).Value = "Lituania"
$ws.Cells.Item(128, 2).Value = 5081
$ws.Cells.Item(128, 3).Value = 125
$ws.Cells.Item(128, 4).Value = 2494
$ws.Cells.Item(128, 5).Value = 2493
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 94

$ws.Cells.Item(129, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(129, 2).Value = 5045
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 4879
$ws.Cells.Item(129, 5).Value = 83
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 83

# Refresh "last updated" timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Octubre de 2020 a las 12:07"
